# Insert a new price-report row at row 282 (weekly update), pushing the
# existing rows 282:368 down to 283:369.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("282:282").Insert()

$ws.Range("A282").Value = 7
$ws.Range("B282").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C282").Value = 'Ñuble'
$ws.Range("D282").Value = 44463
$ws.Range("E282").Value = 16
$ws.Range("F282").Value = 100112033
$ws.Range("G282").Value = 'Lechuga'
$ws.Range("H282").Value = 'Escarola'
$ws.Range("I282").Value = 'Primera'
$ws.Range("J282").Value = 300
$ws.Range("K282").Value = 8500
$ws.Range("L282").Value = 9000
$ws.Range("M282").Value = 8750
$ws.Range("N282").Value = '$/caja 15 unidades'
$ws.Range("O282").Value = 'Provincia del Elquí'
$ws.Range("P282").Value = 583
$ws.Range("Q282").Value = 15
$ws.Range("R282").Value = 'Hortaliza'
